$d = $word.ActiveDocument

# --- Paragraphs 1-2 ("Trace de l'echange..." block): 14pt -> 22pt (sz 28 -> 44) ---
$d.Paragraphs(1).Range.Font.Size = 22
$d.Paragraphs(1).Range.Font.SizeBi = 22
$d.Paragraphs(2).Range.Font.Size = 22
$d.Paragraphs(2).Range.Font.SizeBi = 22

# --- Paragraphs 3-14 (Discord chat transcript block): set to 9pt -> 18pt (sz 36) ---
# Paragraphs that are already non-empty just need their Range font set directly.
# Empty paragraphs need a temporary character inserted so the size "sticks" to the
# paragraph mark's rPr, then that temporary character is removed again.
$emptyIdx = @(3, 6, 9, 12)

For ($i = 3; $i -le 14; $i++) {
    $p = $d.Paragraphs($i)
    if ($emptyIdx -contains $i) {
        $p.Range.InsertBefore("X")
        $p.Range.Font.Size = 18
        $p.Range.Font.SizeBi = 18
        $tmp = $d.Range($p.Range.Start, $p.Range.Start + 1)
        $tmp.Text = ""
    } else {
        $p.Range.Font.Size = 18
        $p.Range.Font.SizeBi = 18
    }
}

# --- Merge the two runs "Un truc dans l'ambiance " + "garage avec un style un peu " ---
# into a single run, by replacing the text unique to the second run with itself; the
# engine coalesces it into the identically-formatted preceding run.
[void]$d.Content.Find.Execute("garage avec un style un peu ", $true, $false, $false, $false, $false, $true, 1, $false, "garage avec un style un peu ", 2)
